$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: P = 16 (Média, already partially present), Q = 17 (new "Situação" column)
$colP = 16
$colQ = 17

# Header row
$ws.Cells.Item(1, $colQ).Value = "Situação"

# Rows whose whole record uses the pre-existing "highlight" row style (cellXfs index 2,
# fill FFFF00 yellow) - on these rows the new Média (P) cell must inherit that same style,
# same as every other cell already in the row.
$customFormatRows = @(5,10,22,29,30,32,36,37,41,42,45,52,53,54,61,68,87,91,96,98,101,103,104,107,112,118,123,126,138,141,143,145,151,152,157)

$green = 9498256   # RGB(144,238,144) -> "Aprovado"
$red   = 8421616   # RGB(240,128,128) -> "Reprovado"

for ($row = 2; $row -le 160; $row++) {
    $pCell = $ws.Cells.Item($row, $colP)
    $pCell.Formula = "=AVERAGE(H$row,N$row,O$row)"

    if ($customFormatRows -contains $row) {
        $rowFillColor = $ws.Cells.Item($row, 1).Interior.Color()
        $pCell.Interior.Color = $rowFillColor
    }

    $avg = $pCell.Value()
    $qCell = $ws.Cells.Item($row, $colQ)

    if ($avg -ge 10) {
        $qCell.Value = "Aprovado"
        $qCell.Interior.Color = $green
    } else {
        $qCell.Value = "Reprovado"
        $qCell.Interior.Color = $red
    }
}

$ws.Range("B2").Select()
